$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values of X1:AA1 (keep their existing style) - Aufgabe 3.2 ergaenzt
$ws.Range("X1:AA1").ClearContents()

# Select AA1 as the active cell
$ws.Range("AA1").Select()
